$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nordex")

# Switch to the Nordex sheet (becomes the active/selected tab on save)
$ws.Activate()

# Update the part number value in B2 (was "WWWW")
$ws.Range("B2").Value = "NORDX-PRT-0001"

# Widen column B to fit the new, longer value
$ws.Columns.Item(2).AutoFit() | Out-Null

# Leave the final selection on D8, as last left by the editor
$ws.Range("D8").Select() | Out-Null
